$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.258.03'
$ws.Range('E2').Value = '  +0.11%  '
$ws.Range('D3').Value = '1.680.56'
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.31'
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5280'
$ws.Range('E6').Value = '  +2.82%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('E8').Value = '  +1.56%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06435'
$ws.Range('E9').Value = '  +0.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.05'
$ws.Range('E10').Value = '  +2.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07498'
$ws.Range('E11').Value = '  +1.55%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.548'
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.676.63'
$ws.Range('E13').Value = '  +0.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5810'
$ws.Range('E14').Value = '  -0.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008462'
$ws.Range('E15').Value = '  -2.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.31'
$ws.Range('E16').Value = '  -0.63%  '
$ws.Range('D17').Value = '26.322.78'
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.920'
$ws.Range('E18').Value = '  -0.92%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.88'
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '189.42'
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.204'
$ws.Range('E22').Value = '  -0.27%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.009'
$ws.Range('E23').Value = '  +0.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '144.95'
$ws.Range('E24').Value = '  +0.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.722'
$ws.Range('E25').Value = '  +1.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1236'
$ws.Range('E26').Value = '  +4.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.80'
$ws.Range('E27').Value = '  +1.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06597'
$ws.Range('E28').Value = '  +10.54%  '
$ws.Range('E29').Value = '  +5.83%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.328'
$ws.Range('E30').Value = '  +0.52%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.583'
$ws.Range('E31').Value = '  +1.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.570'
$ws.Range('E32').Value = '  +1.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.661'
$ws.Range('E33').Value = '  +1.40%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.025'
$ws.Range('E34').Value = '  +0.95%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6198'
$ws.Range('E35').Value = '  +2.86%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.398'
$ws.Range('E36').Value = '  +1.00%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.699'
$ws.Range('E37').Value = '  +1.94%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.392'
$ws.Range('E38').Value = '  +4.86%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01622'
$ws.Range('E39').Value = '  +0.11%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '1.107.01'
$ws.Range('E40').Value = '  +2.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8771'
$ws.Range('E41').Value = '  +0.79%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.015'
$ws.Range('E42').Value = '  +0.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.54'
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('D44').Value = '1.827.92'
$ws.Range('E44').Value = '  +0.35%  '
$ws.Range('E45').Value = '  -2.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.85'
$ws.Range('E46').Value = '  +1.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.167'
$ws.Range('E47').Value = '  +1.50%  '
$ws.Range('E48').Value = '  -0.39%  '
$ws.Range('E50').Value = '  +0.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.042'
$ws.Range('E51').Value = '  +2.53%  '
